# Power_Network.xlsx — switch the line "Technical Representation" values
# from the DC-OPF approximation to SOCP (second-order cone), per the
# commit "Add quadratic constraint (probably bugged) and quadratic obj
# handeling".
#
# Column O ("pTecRepr" / "Technical Representation") on sheet ScenarioA
# holds one literal value per network line (rows 8-20); all of them read
# "DC-OPF" and need to become "SOCP".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenarioA")

$ws.Range("O8:O20").Value = "SOCP"

# Leave the cursor roughly where the author ended up after editing the
# column (just below the data table).
$ws.Range("O26").Select() | Out-Null
